$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures and fix two mis-ordered rows
# (WrappedBTC/WrappedEther swap at rows 18-19, VeChain/Kaspa swap at rows 41-42).
# Numeric-looking price strings are forced to Text (NumberFormat "@" then
# reset to the default Normal style) so they keep their original literal
# formatting (leading zeros / trailing zeros / thousands-dot grouping)
# instead of being auto-converted to floating point numbers by Excel.

$ws.Range("D2").Value = "67.737.24"
$ws.Range("E2").Value = "  -4.08%  "
$ws.Range("D3").Value = "3.440.63"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.601"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.91%  "
$ws.Range("D8").Value = "3.433.76"
$ws.Range("E8").Value = "  -5.34%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  -8.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000269"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.56%  "
$ws.Range("D15").Value = "3.994.06"
$ws.Range("E15").Value = "  -5.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "610.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.87%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.460.60"
$ws.Range("E18").Value = "  -4.68%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "67.810.22"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.865"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "94.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.28%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  -7.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.33%  "
$ws.Range("E32").Value = "  -8.31%  "
$ws.Range("E33").Value = "  -7.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "588.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("E36").Value = "  -4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.15%  "
$ws.Range("E40").Value = "  -14.98%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.134"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.24%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0428"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.96%  "
$ws.Range("D43").Value = "3.357.07"
$ws.Range("E43").Value = "  -5.04%  "
$ws.Range("E44").Value = "  -8.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.87%  "
$ws.Range("D46").Value = "0.0₃0675"
$ws.Range("E46").Value = "  -7.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.97%  "
$ws.Range("E49").Value = "  -6.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.84%  "
